$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D is column index 4
$col = 4

$updates = @{
    2  = 6556.266582378041
    3  = 6556.266582378034
    4  = 6556.266582378034
    6  = 18732.19023536615
    7  = 18732.19023536615
    9  = 2062.858975181512
    10 = 2062.85897518151
    11 = 30
    12 = 1702.85897518151
    13 = 30
    17 = 24000
    19 = 1485.548159853576
    20 = 1485.548159853576
    21 = 360
    24 = 283211.0373750483
    25 = 283211.037375048
    28 = 5664.220747500913
    29 = 5664.220747500913
    30 = 283211.037375048
    35 = 20217.60000000103
    36 = 20217.60000000103
    38 = 0
    39 = 0
    41 = 0
    42 = 0
    43 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, $col).Value = $updates[$row]
}
